# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled
# update). Every assigned value is prefixed with a literal leading
# apostrophe (PowerShell single-quoted string, so '' -> a literal ')
# to force Excel to keep numeric-looking text (e.g. "312.55", "41.469.60")
# as plain text instead of auto-coercing it to a number, exactly like the
# source workbook stores these as inline/shared strings. The Style reset
# to 'Normal' right after strips the transient quote-prefix formatting
# Excel applies for that entry mode, so the cell keeps its original
# (unstyled) appearance.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''41.469.60'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  -3.09%  '
$ws.Range('E2').Style = 'Normal'

$ws.Range('D3').Value = '''2.475.95'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  -2.70%  '
$ws.Range('E3').Style = 'Normal'

$ws.Range('E4').Value = '''  +0.20%  '
$ws.Range('E4').Style = 'Normal'

$ws.Range('D5').Value = '''312.55'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  +0.07%  '
$ws.Range('E5').Style = 'Normal'

$ws.Range('D6').Value = '''94.43'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '''  -6.30%  '
$ws.Range('E6').Style = 'Normal'

$ws.Range('E7').Value = '''  -3.20%  '
$ws.Range('E7').Style = 'Normal'

$ws.Range('E8').Value = '''  +0.20%  '
$ws.Range('E8').Style = 'Normal'

$ws.Range('D9').Value = '''0.499'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '''  -4.49%  '
$ws.Range('E9').Style = 'Normal'

$ws.Range('D10').Value = '''33.54'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''  -5.59%  '
$ws.Range('E10').Style = 'Normal'

$ws.Range('E11').Value = '''  -2.79%  '
$ws.Range('E11').Style = 'Normal'

$ws.Range('E12').Value = '''  -0.90%  '
$ws.Range('E12').Style = 'Normal'

$ws.Range('D13').Value = '''6.99'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''  -4.37%  '
$ws.Range('E13').Style = 'Normal'

$ws.Range('D14').Value = '''2.860.44'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''  -2.57%  '
$ws.Range('E14').Style = 'Normal'

$ws.Range('B15').Value = '''Chainlink'
$ws.Range('B15').Style = 'Normal'
$ws.Range('C15').Value = '''https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('C15').Style = 'Normal'
$ws.Range('D15').Value = '''15.29'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  -0.75%  '
$ws.Range('E15').Style = 'Normal'

$ws.Range('B16').Value = '''WrappedEther'
$ws.Range('B16').Style = 'Normal'
$ws.Range('C16').Value = '''https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('C16').Style = 'Normal'
$ws.Range('D16').Value = '''2.427.07'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  -4.55%  '
$ws.Range('E16').Style = 'Normal'

$ws.Range('D17').Value = '''0.789'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '''  -3.53%  '
$ws.Range('E17').Style = 'Normal'

$ws.Range('D18').Value = '''41.401.30'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  -3.24%  '
$ws.Range('E18').Style = 'Normal'

$ws.Range('D19').Value = '''6.32'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '''  -6.32%  '
$ws.Range('E19').Style = 'Normal'

$ws.Range('D20').Value = '''0.0₃0926'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  -2.71%  '
$ws.Range('E20').Style = 'Normal'

$ws.Range('D21').Value = '''11.25'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''  -8.79%  '
$ws.Range('E21').Style = 'Normal'

$ws.Range('D22').Value = '''68.77'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '''  -1.71%  '
$ws.Range('E22').Style = 'Normal'

$ws.Range('D23').Value = '''237.39'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '''  -2.41%  '
$ws.Range('E23').Style = 'Normal'

$ws.Range('D24').Value = '''2.75'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''  -4.49%  '
$ws.Range('E24').Style = 'Normal'

$ws.Range('D26').Value = '''1.90'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '''  -6.09%  '
$ws.Range('E26').Style = 'Normal'

$ws.Range('D27').Value = '''24.08'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '''  -6.20%  '
$ws.Range('E27').Style = 'Normal'

$ws.Range('D28').Value = '''2.25'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '''  -3.88%  '
$ws.Range('E28').Style = 'Normal'

$ws.Range('D29').Value = '''9.70'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '''  -4.54%  '
$ws.Range('E29').Style = 'Normal'

$ws.Range('D30').Value = '''36.46'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '''  -5.16%  '
$ws.Range('E30').Style = 'Normal'

$ws.Range('D31').Value = '''152.18'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '''  -3.78%  '
$ws.Range('E31').Style = 'Normal'

$ws.Range('D32').Value = '''5.47'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '''  -7.04%  '
$ws.Range('E32').Style = 'Normal'

$ws.Range('D33').Value = '''2.57'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '''  -3.72%  '
$ws.Range('E33').Style = 'Normal'

$ws.Range('E34').Value = '''  -6.96%  '
$ws.Range('E34').Style = 'Normal'

$ws.Range('D35').Value = '''0.0747'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '''  -5.99%  '
$ws.Range('E35').Style = 'Normal'

$ws.Range('E36').Value = '''  -2.59%  '
$ws.Range('E36').Style = 'Normal'

$ws.Range('D37').Value = '''17.49'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '''  -2.54%  '
$ws.Range('E37').Style = 'Normal'

$ws.Range('E38').Value = '''  -5.23%  '
$ws.Range('E38').Style = 'Normal'

$ws.Range('E39').Value = '''  -2.82%  '
$ws.Range('E39').Style = 'Normal'

$ws.Range('B40').Value = '''Kaspa'
$ws.Range('B40').Style = 'Normal'
$ws.Range('C40').Value = '''https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('C40').Style = 'Normal'
$ws.Range('D40').Value = '''0.101'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  -8.70%  '
$ws.Range('E40').Style = 'Normal'

$ws.Range('B41').Value = '''RenderToken'
$ws.Range('B41').Style = 'Normal'
$ws.Range('C41').Value = '''https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('C41').Style = 'Normal'
$ws.Range('D41').Value = '''4.24'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  +2.28%  '
$ws.Range('E41').Style = 'Normal'

$ws.Range('D42').Value = '''1.01'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '''  +0.44%  '
$ws.Range('E42').Style = 'Normal'

$ws.Range('D43').Value = '''19.49'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''  -11.23%  '
$ws.Range('E43').Style = 'Normal'

$ws.Range('D44').Value = '''1.986.83'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '''  -0.64%  '
$ws.Range('E44').Style = 'Normal'

$ws.Range('E45').Value = '''  -4.41%  '
$ws.Range('E45').Style = 'Normal'

$ws.Range('D46').Value = '''3.00'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''  -8.62%  '
$ws.Range('E46').Style = 'Normal'

$ws.Range('D47').Value = '''8.74'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '''  -4.58%  '
$ws.Range('E47').Style = 'Normal'

$ws.Range('D48').Value = '''2.723.40'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '''  -2.17%  '
$ws.Range('E48').Style = 'Normal'

$ws.Range('D49').Value = '''69.56'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '''  -4.05%  '
$ws.Range('E49').Style = 'Normal'

$ws.Range('D50').Value = '''97.16'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '''  -4.36%  '
$ws.Range('E50').Style = 'Normal'

$ws.Range('D51').Value = '''74.55'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '''  -6.74%  '
$ws.Range('E51').Style = 'Normal'
